$d = $word.ActiveDocument

# Step 1: split the run " e trabalho para nota" into three runs with a proofErr-wrapped "trabalho"
$para2 = $d.Paragraphs.Item(2)
$xmlPara2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>Projeto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>desenvolvido</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>como</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>desafio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>trabalho</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> para nota</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.Range.InsertXML($xmlPara2)

# Step 2: insert 6 empty paragraphs before the existing "Testes" heading paragraph (paragraph 3)
$anchor = $d.Paragraphs.Item(3).Range
$anchor.Collapse(1)
for ($i = 0; $i -lt 6; $i++) {
    $anchor.InsertParagraphBefore()
}

# Step 3: fill each newly created paragraph with its target content
$xmlNew0 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:spacing w:before="0" w:beforeAutospacing="off" w:after="160" w:afterAutospacing="off" w:line="259" w:lineRule="auto"/><w:ind w:left="0" w:right="0"/><w:jc w:val="left"/><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Testes</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(3).Range.InsertXML($xmlNew0)
$xmlNew1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Marcelo Araujo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(4).Range.InsertXML($xmlNew1)
$xmlNew2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr/><w:t>Claiver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr/><w:t xml:space="preserve"> Bina</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(5).Range.InsertXML($xmlNew2)
$xmlNew3 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Pedro Henrique</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(6).Range.InsertXML($xmlNew3)
$xmlNew4 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Bruno Zuppa</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(7).Range.InsertXML($xmlNew4)
$xmlNew5 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Leandro Henrique de Lima</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item(8).Range.InsertXML($xmlNew5)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text
}
